# Generate Report for Handback
# -----------------------------------------------------------------------
# This script mirrors the xlsx-level diff:
#  1. Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#     on Overview!E2:F2/E3:F3 and zh-cn!C2:C3 / de-de!C2:C3.
#  2. zh-cn "Latest Handback DateTime" (K2:K3) goes from the zero date to
#     "2016-08-25 22:36:47"; de-de's goes to "2016-08-25 22:36:54".
#  3. Both language sheets gain a "Latest Target File" (I) / "Latest
#     Handback File" (J) entry: I = "a.md" hyperlinked back to the source
#     markdown file (same target as the existing A-column link), J = the
#     generated handback xliff file name.
#  4. Column widths widen to fit the new, longer text.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cf504e351a46d87b9eedbafcc69242edaca656c0/e2e/a.md"

# ---------------------------------------------------------------------
# Overview sheet: just the status text refresh.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# Helper applied identically to the zh-cn and de-de sheets.
# ---------------------------------------------------------------------
function Set-HandbackSheet($SheetName, $HandbackDateTime, $XliffFileName) {

    $ws = $wb.Worksheets.Item($SheetName)

    # Status column
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Latest Target File (I) : hyperlinked "a.md" pointing at the source doc
    $ws.Range("I2").Value = "a.md"
    $ws.Range("I3").Value = "a.md"
    $ws.Range("I2:I3").Font.Underline = $true
    $ws.Range("I2:I3").Font.Color = 15570276
    $ws.Hyperlinks.Add($ws.Range("I2"), $baseUrl, [Type]::Missing, [Type]::Missing, "a.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $baseUrl, [Type]::Missing, [Type]::Missing, "a.md")

    # Latest Handback File (J) : generated xliff file name
    $ws.Range("J2").Value = $XliffFileName
    $ws.Range("J3").Value = $XliffFileName

    # Latest Handback DateTime (K)
    $ws.Range("K2").Value = $HandbackDateTime
    $ws.Range("K3").Value = $HandbackDateTime

    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}

Set-HandbackSheet "zh-cn" "2016-08-25 22:36:47" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
Set-HandbackSheet "de-de" "2016-08-25 22:36:54" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
